{"js": "/*\n * Replace the math-expression text in each cell of the single table,\n * in document order (row-major), matching the values from the diff.\n * Each cell has exactly one paragraph with one run; we only rewrite the\n * run's text (via paragraph.insertText/replace) so that the run/paragraph\n * formatting (rFonts, sz, jc) is preserved untouched, matching the diff\n * which only touches the <w:t> contents.\n */\nconst oldValues = [\"44-28=\", \"35+58=\", \"37+59=\", \"62-13=\", \"78-45=\", \"84-80=\", \"78+12=\", \"64+0=\", \"77+22=\", \"34+47=\", \"46-13=\", \"28-4=\", \"4+81=\", \"83+9=\", \"9+20=\", \"28+2=\", \"51-6=\", \"34-31=\", \"1+80=\", \"72-45=\", \"33+26=\", \"41+16=\", \"6+57=\", \"71+17=\", \"18+59=\", \"65+22=\", \"32+38=\", \"36+7=\", \"9+18=\", \"51-28=\", \"67-60=\", \"19+44=\", \"71+16=\", \"18+71=\", \"6+47=\", \"80-32=\", \"90-51=\", \"30+40=\", \"52+16=\", \"70-66=\", \"97-8=\", \"13+5=\", \"56-49=\", \"82-31=\", \"89-18=\", \"18-10=\", \"95-61=\", \"59-19=\", \"68+1=\", \"2+44=\", \"77-42=\", \"61-23=\", \"46+13=\", \"20-4=\", \"81-51=\", \"31+14=\", \"11+41=\", \"15+56=\", \"7+73=\", \"13+83=\", \"25+37=\", \"27+57=\", \"55-8=\", \"6+43=\", \"61-55=\", \"98-90=\", \"16+31=\", \"43-26=\", \"89+2=\", \"59+34=\", \"14-4=\", \"60-16=\", \"9-4=\", \"12+19=\", \"19-8=\", \"20-10=\", \"8+89=\", \"92-6=\", \"12+66=\", \"6+9=\", \"4+73=\", \"73+6=\", \"11+64=\", \"69-13=\", \"86-39=\", \"94-83=\", \"25-7=\", \"11+50=\", \"15+50=\", \"95-76=\", \"59-56=\", \"29+8=\", \"12-7=\", \"88-14=\", \"4+86=\", \"58+0=\", \"8-5=\", \"70-44=\", \"77-74=\", \"57+14=\"];\nconst newValues = [\"37-30=\", \"15+2=\", \"40+12=\", \"25+3=\", \"37-34=\", \"15+40=\", \"44+37=\", \"87-64=\", \"40+29=\", \"12+18=\", \"70+17=\", \"26+33=\", \"16+76=\", \"67-45=\", \"47+5=\", \"4+38=\", \"71-21=\", \"51-42=\", \"76-54=\", \"31-27=\", \"14+6=\", \"49-28=\", \"16+71=\", \"6+78=\", \"52+11=\", \"49+16=\", \"22+58=\", \"35+35=\", \"52-50=\", \"43+48=\", \"44+48=\", \"32+64=\", \"40+37=\", \"59-4=\", \"52-25=\", \"83-68=\", \"59-32=\", \"15-0=\", \"45-17=\", \"21+11=\", \"54-28=\", \"66-7=\", \"48-3=\", \"42+51=\", \"71-30=\", \"2+83=\", \"75-17=\", \"6+74=\", \"12+15=\", \"86-56=\", \"92-9=\", \"17+19=\", \"63-8=\", \"70-2=\", \"42-1=\", \"24-10=\", \"41-37=\", \"7+15=\", \"58-28=\", \"21+9=\", \"71-58=\", \"53-10=\", \"60-31=\", \"12+57=\", \"61-27=\", \"88-18=\", \"7+48=\", \"65+5=\", \"4-4=\", \"11+40=\", \"38+60=\", \"30+62=\", \"14+55=\", \"6+93=\", \"86-59=\", \"34-11=\", \"29+28=\", \"45+20=\", \"83-64=\", \"12+38=\", \"22+58=\", \"5+77=\", \"20+20=\", \"56+18=\", \"32+20=\", \"11+80=\", \"76-27=\", \"39+48=\", \"45-33=\", \"16+4=\", \"35+25=\", \"54-22=\", \"46-15=\", \"65-35=\", \"13+84=\", \"87-61=\", \"14+82=\", \"23+76=\", \"2+29=\", \"92+3=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document\");\n}\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rows = table.rowCount;\nconst cols = rows > 0 ? table.values[0].length : 0;\n\n// Collect all paragraphs (one per cell) first so we can batch-load their text.\nconst paragraphs = [];\nfor (let r = 0; r < rows; r++) {\n  for (let c = 0; c < cols; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.load(\"text\");\n    paragraphs.push(para);\n  }\n}\nawait context.sync();\n\nif (paragraphs.length !== newValues.length) {\n  throw new Error(\n    \"Cell count (\" + paragraphs.length + \") does not match expected (\" + newValues.length + \")\"\n  );\n}\n\nfor (let i = 0; i < paragraphs.length; i++) {\n  const para = paragraphs[i];\n  const current = para.text;\n  const expectedOld = oldValues[i];\n  const newValue = newValues[i];\n  // Only touch paragraphs whose current text still matches the expected\n  // \"before\" value -- this keeps the edit targeted/idempotent even if\n  // run order differs slightly from assumption.\n  if (current === newValue) {\n    continue; // already updated\n  }\n  if (current !== expectedOld) {\n    throw new Error(\n      \"Cell \" + i + \" text mismatch: expected '\" + expectedOld + \"' but found '\" + current + \"'\"\n    );\n  }\n  para.insertText(newValue, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the math-expression text in each cell of the single table,\n# in document order (row-major), matching the values from the diff.\n# Each cell's Range.Text is reassigned directly (Word COM keeps the\n# existing run/paragraph formatting -- rFonts, sz, jc -- and only the\n# <w:t> text content changes), mirroring the diff which only touches\n# the text nodes.\n\n$oldValues = @(\"44-28=\", \"35+58=\", \"37+59=\", \"62-13=\", \"78-45=\", \"84-80=\", \"78+12=\", \"64+0=\", \"77+22=\", \"34+47=\", \"46-13=\", \"28-4=\", \"4+81=\", \"83+9=\", \"9+20=\", \"28+2=\", \"51-6=\", \"34-31=\", \"1+80=\", \"72-45=\", \"33+26=\", \"41+16=\", \"6+57=\", \"71+17=\", \"18+59=\", \"65+22=\", \"32+38=\", \"36+7=\", \"9+18=\", \"51-28=\", \"67-60=\", \"19+44=\", \"71+16=\", \"18+71=\", \"6+47=\", \"80-32=\", \"90-51=\", \"30+40=\", \"52+16=\", \"70-66=\", \"97-8=\", \"13+5=\", \"56-49=\", \"82-31=\", \"89-18=\", \"18-10=\", \"95-61=\", \"59-19=\", \"68+1=\", \"2+44=\", \"77-42=\", \"61-23=\", \"46+13=\", \"20-4=\", \"81-51=\", \"31+14=\", \"11+41=\", \"15+56=\", \"7+73=\", \"13+83=\", \"25+37=\", \"27+57=\", \"55-8=\", \"6+43=\", \"61-55=\", \"98-90=\", \"16+31=\", \"43-26=\", \"89+2=\", \"59+34=\", \"14-4=\", \"60-16=\", \"9-4=\", \"12+19=\", \"19-8=\", \"20-10=\", \"8+89=\", \"92-6=\", \"12+66=\", \"6+9=\", \"4+73=\", \"73+6=\", \"11+64=\", \"69-13=\", \"86-39=\", \"94-83=\", \"25-7=\", \"11+50=\", \"15+50=\", \"95-76=\", \"59-56=\", \"29+8=\", \"12-7=\", \"88-14=\", \"4+86=\", \"58+0=\", \"8-5=\", \"70-44=\", \"77-74=\", \"57+14=\")\n$newValues = @(\"37-30=\", \"15+2=\", \"40+12=\", \"25+3=\", \"37-34=\", \"15+40=\", \"44+37=\", \"87-64=\", \"40+29=\", \"12+18=\", \"70+17=\", \"26+33=\", \"16+76=\", \"67-45=\", \"47+5=\", \"4+38=\", \"71-21=\", \"51-42=\", \"76-54=\", \"31-27=\", \"14+6=\", \"49-28=\", \"16+71=\", \"6+78=\", \"52+11=\", \"49+16=\", \"22+58=\", \"35+35=\", \"52-50=\", \"43+48=\", \"44+48=\", \"32+64=\", \"40+37=\", \"59-4=\", \"52-25=\", \"83-68=\", \"59-32=\", \"15-0=\", \"45-17=\", \"21+11=\", \"54-28=\", \"66-7=\", \"48-3=\", \"42+51=\", \"71-30=\", \"2+83=\", \"75-17=\", \"6+74=\", \"12+15=\", \"86-56=\", \"92-9=\", \"17+19=\", \"63-8=\", \"70-2=\", \"42-1=\", \"24-10=\", \"41-37=\", \"7+15=\", \"58-28=\", \"21+9=\", \"71-58=\", \"53-10=\", \"60-31=\", \"12+57=\", \"61-27=\", \"88-18=\", \"7+48=\", \"65+5=\", \"4-4=\", \"11+40=\", \"38+60=\", \"30+62=\", \"14+55=\", \"6+93=\", \"86-59=\", \"34-11=\", \"29+28=\", \"45+20=\", \"83-64=\", \"12+38=\", \"22+58=\", \"5+77=\", \"20+20=\", \"56+18=\", \"32+20=\", \"11+80=\", \"76-27=\", \"39+48=\", \"45-33=\", \"16+4=\", \"35+25=\", \"54-22=\", \"46-15=\", \"65-35=\", \"13+84=\", \"87-61=\", \"14+82=\", \"23+76=\", \"2+29=\", \"92+3=\")\n\n$d = $word.ActiveDocument\n\nif ($d.Tables.Count -lt 1) {\n    throw \"No table found in document\"\n}\n\n$tbl = $d.Tables.Item(1)\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\nif (($rows * $cols) -ne $newValues.Length) {\n    throw \"Cell count ($($rows * $cols)) does not match expected ($($newValues.Length))\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $current = $cell.Range.Text.TrimEnd([char]7, [char]13)\n        $expectedOld = $oldValues[$i]\n        $newValue = $newValues[$i]\n\n        if ($current -eq $newValue) {\n            # already updated\n            $i++\n            continue\n        }\n        if ($current -ne $expectedOld) {\n            throw \"Cell $i (row $r, col $c) text mismatch: expected '$expectedOld' but found '$current'\"\n        }\n        $cell.Range.Text = $newValue\n        $i++\n    }\n}\n"}
